$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields (row 2-3) ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long, purely-numeric-looking card number that must stay a TEXT
# cell (like the original). A plain .Value assignment would be stored as a
# number, so write it through a text formula and then collapse the formula
# down to a static value via copy/paste-values (keeps style s=8 intact and
# avoids forcing a NumberFormat="@"/quote-prefix style switch).
$ws.Range("B3").Formula = '="2570314725427075"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 06.01.2025"

# --- Row 6 ---
$ws.Range("B6").Value = "07.01."
$ws.Range("C6").Value = "08.01."
$ws.Range("D6").Value = "BURGER KING Aschaffenburg"
$ws.Range("E6").Value = "8,84-"

# --- Row 7 ---
$ws.Range("B7").Value = "10.01."
$ws.Range("C7").Value = "11.01."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 87231101"
$ws.Range("E7").Value = "86,37-"

# --- Row 8 ---
$ws.Range("B8").Value = "11.01."
$ws.Range("C8").Value = "12.01."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 87701066"
$ws.Range("E8").Value = "38,10-"

# --- Row 9 (previously a blank filler row; becomes a real transaction) ---
$ws.Range("B9").Value = "13.01."
$ws.Range("C9").Value = "14.01."
$ws.Range("D9").Value = "PAYPAL GCCOUY"
# E9 needs to switch from the blank-row style (s=13) to the amount-column
# style (s=17, right-aligned) used by the other populated rows (E6-E8).
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E9").Value = "98,09-"

# --- Row 10 (previously a blank filler row; becomes a real transaction) ---
$ws.Range("B10").Value = "14.01."
$ws.Range("C10").Value = "15.01."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
# E10 needs to switch from the blank-row style (s=12) to the amount-column
# style (s=17, right-aligned) used by the other populated rows (E6-E8).
$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E10").Value = "24,70-"

$excel.CutCopyMode = $false

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 18.01.2025"
$ws.Range("E12").Value = "256,10-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 27.01.2025"
